# "arrange checkout page and comments for save details and sort products"
# - Replace the row-2 product (code/name) with the new product entry.
# - Remove the old rows 3 and 4 (the sheet now only lists a single product row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new product code/name (leading apostrophe keeps the zero-padded,
# space-padded code stored as text rather than being coerced to a number).
$ws.Range("A2").Value = "'0010017000044       "
$ws.Range("B2").Value = "ABRAZADERA MARCO SILLA POLISPORT                  "
$ws.Range("C2").Value = 2

# Drop the two trailing product rows entirely (dimension shrinks to A1:C2).
$ws.Rows("3:4").Delete()
